# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the cell to carry an explicit text value (matches the source data,
    # which stores every Price/Volume cell as a string, even the numeric-looking ones)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '66.965.15'
$ws.Range('E2').Value = '  +5.01%  '

Set-TextValue $ws.Range('D3') '3.509.12'
$ws.Range('E3').Value = '  +2.78%  '

$ws.Range('E4').Value = '  +0.01%  '

Set-TextValue $ws.Range('D5') '593.76'
$ws.Range('E5').Value = '  +4.16%  '

Set-TextValue $ws.Range('D6') '169.01'
$ws.Range('E6').Value = '  +7.26%  '

Set-TextValue $ws.Range('D7') '0.999'
$ws.Range('E7').Value = '  -0.07%  '

Set-TextValue $ws.Range('D8') '3.507.95'
$ws.Range('E8').Value = '  +2.67%  '

Set-TextValue $ws.Range('D9') '0.574'
$ws.Range('E9').Value = '  +1.17%  '

Set-TextValue $ws.Range('D10') '7.28'
$ws.Range('E10').Value = '  +0.74%  '

$ws.Range('E11').Value = '  +5.41%  '

Set-TextValue $ws.Range('D12') '0.440'
$ws.Range('E12').Value = '  +4.41%  '

Set-TextValue $ws.Range('D13') '4.112.68'
$ws.Range('E13').Value = '  +2.80%  '

Set-TextValue $ws.Range('D14') '0.135'
$ws.Range('E14').Value = '  +0.17%  '

Set-TextValue $ws.Range('D15') '28.27'
$ws.Range('E15').Value = '  +4.37%  '

$ws.Range('E16').Value = '  +4.48%  '

Set-TextValue $ws.Range('D17') '66.901.32'
$ws.Range('E17').Value = '  +4.77%  '

Set-TextValue $ws.Range('D18') '3.498.34'
$ws.Range('E18').Value = '  +2.26%  '

Set-TextValue $ws.Range('D19') '6.33'
$ws.Range('E19').Value = '  +4.25%  '

Set-TextValue $ws.Range('D20') '14.06'
$ws.Range('E20').Value = '  +3.10%  '

Set-TextValue $ws.Range('D21') '395.39'
$ws.Range('E21').Value = '  +3.06%  '

Set-TextValue $ws.Range('D22') '7.97'
$ws.Range('E22').Value = '  +2.36%  '

Set-TextValue $ws.Range('D23') '73.16'
$ws.Range('E23').Value = '  +2.60%  '

Set-TextValue $ws.Range('D24') '0.0000128'
$ws.Range('E24').Value = '  +11.61%  '

$ws.Range('E25').Value = '  -0.20%  '

$ws.Range('E26').Value = '  +3.08%  '

Set-TextValue $ws.Range('D27') '10.14'
$ws.Range('E27').Value = '  +5.15%  '

$ws.Range('E28').Value = '  +2.37%  '

$ws.Range('E29').Value = '  +0.10%  '

Set-TextValue $ws.Range('D30') '6.42'
$ws.Range('E30').Value = '  +5.25%  '

$ws.Range('E31').Value = '  +6.26%  '

$ws.Range('E32').Value = '  +4.18%  '

$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range('D33') '23.60'
$ws.Range('E33').Value = '  +3.09%  '

$ws.Range('B34').Value = 'Aptos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D34') '7.46'
$ws.Range('E34').Value = '  +7.43%  '

Set-TextValue $ws.Range('D35') '0.999'
$ws.Range('E35').Value = '  +0.05%  '

$ws.Range('E36').Value = '  +6.75%  '

Set-TextValue $ws.Range('D37') '162.06'
$ws.Range('E37').Value = '  +0.82%  '

$ws.Range('E38').Value = '  +6.27%  '

$ws.Range('E39').Value = '  +6.53%  '

$ws.Range('E40').Value = '  +3.96%  '

Set-TextValue $ws.Range('D41') '4.67'
$ws.Range('E41').Value = '  +7.31%  '

Set-TextValue $ws.Range('D42') '2.850.83'
$ws.Range('E42').Value = '  +1.93%  '

Set-TextValue $ws.Range('D43') '26.48'
$ws.Range('E43').Value = '  +2.17%  '

Set-TextValue $ws.Range('D44') '6.70'
$ws.Range('E44').Value = '  +4.91%  '

Set-TextValue $ws.Range('D45') '43.47'
$ws.Range('E45').Value = '  +1.11%  '

Set-TextValue $ws.Range('D46') '26.49'
$ws.Range('E46').Value = '  +1.11%  '

$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range('D47') '2.56'
$ws.Range('E47').Value = '  +6.23%  '

$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D48') '0.0316'
$ws.Range('E48').Value = '  +4.25%  '

Set-TextValue $ws.Range('D49') '350.51'
$ws.Range('E49').Value = '  +6.14%  '

Set-TextValue $ws.Range('D50') '1.09'
$ws.Range('E50').Value = '  +5.26%  '

Set-TextValue $ws.Range('D51') '33.70'
$ws.Range('E51').Value = '  +12.37%  '
